$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the "ueiDUNS" column (column C) and remove it entirely, shifting
# subsequent columns left (this also drops the now-unused "ueiDUNS" shared string)
$col = $ws.Columns.Item(3)
$col.Select() | Out-Null
$col.Delete() | Out-Null
